$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the ")" paragraph that closes the "CREATE TABLE composition_equipe"
# statement (it is the paragraph right after the
# "FOREIGN KEY (idUtilisateur) REFERENCES utilisateurs(idUtilisateur)" line),
# and the blank paragraph that follows it.
# ---------------------------------------------------------------------------
$fkParagraph = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*FOREIGN KEY (idUtilisateur) REFERENCES utilisateurs(idUtilisateur)*") {
        $fkParagraph = $p
    }
}

$closeParen = $fkParagraph.Next()
$blankAfter = $closeParen.Next()

# Make room for the new content by inserting an empty paragraph right before
# the existing blank line, then fill the whole new block (closing ")"
# duplicate + the new participants_evenements table + the trailing blank
# line) in one shot using CR-separated text; Word splits it into individual
# paragraphs for us.
$blankAfter.Range.InsertParagraphBefore()

$insertPoint = $d.Range($blankAfter.Range.Start, $blankAfter.Range.Start)

$block = ")`rCREATE TABLE participants_evenements`r(`r `tidUtilisateur int NOT NULL,`r    idEvenement int NOT NULL,`r    FOREIGN KEY (idUtilisateur) REFERENCES utilisateurs(idUtilisateur),`r    FOREIGN KEY (idEvenement) REFERENCES evenements(id)`r)`r"

$insertPoint.InsertBefore($block)
